# Update gh-pages to output generated at 456a3b4
# Applies numeric "想去人数" (F column) bumps and a couple of "最低票价"
# (G column) status-text changes ("已售罄" -> "不可售") across the
# 展览 / 演出 / 本地生活 / 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

$wsExpo   = $wb.Worksheets.Item("展览")
$wsShow   = $wb.Worksheets.Item("演出")
$wsLocal  = $wb.Worksheets.Item("本地生活")
$wsAll    = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExpo.Range("F3").Value = 26753
$wsExpo.Range("F4").Value = 591
$wsExpo.Range("F5").Value = 256
$wsExpo.Range("F6").Value = 615
$wsExpo.Range("F7").Value = 177
$wsExpo.Range("F8").Value = 554
$wsExpo.Range("F11").Value = 244
$wsExpo.Range("F12").Value = 191
$wsExpo.Range("F15").Value = 72
$wsExpo.Range("F16").Value = 432
$wsExpo.Range("F18").Value = 1552
$wsExpo.Range("F19").Value = 212

# --- 演出 (sheet2) ---
$wsShow.Range("G2").Value = "不可售"
$wsShow.Range("F3").Value = 234

# --- 本地生活 (sheet3) ---
$wsLocal.Range("F2").Value = 5087

# --- 全部类型 (sheet4) ---
$wsAll.Range("F3").Value = 5087
$wsAll.Range("F5").Value = 26753
$wsAll.Range("F6").Value = 591
$wsAll.Range("G7").Value = "不可售"
$wsAll.Range("F8").Value = 256
$wsAll.Range("F9").Value = 234
$wsAll.Range("F10").Value = 615
$wsAll.Range("F13").Value = 177
$wsAll.Range("F19").Value = 554
$wsAll.Range("F23").Value = 244
$wsAll.Range("F24").Value = 191
$wsAll.Range("F28").Value = 72
$wsAll.Range("F31").Value = 432
$wsAll.Range("F34").Value = 1552
$wsAll.Range("F35").Value = 212
